# Applies the diff: append a new "4:" entry to the Overblik.docx work log.
#
# The document currently ends with:
#   ... "3:" paragraph ...
#   ... elastix paragraph ...
#   <empty paragraph>
#   <sectPr>
#
# We need to add, after that trailing empty paragraph and before the
# section properties, two new paragraphs:
#   "4:"
#   "Croppede en lille fraktur ud. ... lort uanset hvad." + " Tror jeg prøver med et større volumen nu.."
# where the second paragraph is made of two separate runs (the second
# run carries a leading space and must render with xml:space="preserve").

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# Grab the range right at the end of the document's last (empty) paragraph
# and use it to create two brand-new, empty paragraphs after it without
# disturbing that existing empty paragraph.
$tailRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()
$tailRange.InsertParagraphAfter()

# The two freshly-created paragraphs are now the last two in the document.
$p4 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$p5 = $d.Paragraphs.Item($d.Paragraphs.Count)

# First new paragraph: "4:"
$xml4 = "<w:p xmlns:w='$wNs'><w:r><w:t>4:</w:t></w:r></w:p>"
$p4.Range.InsertXML($xml4)

# Second new paragraph: two runs - the long note, then a trailing
# " Tror jeg prøver med et større volumen nu.." run with a preserved
# leading space.
$run1 = "Croppede en lille fraktur ud. Med samme parameterfile som normalt. Uden at normalisere. Og det bliver virkelig lort. Den har svært ved at gøre det ordentligt. Prøvede at normalisere til [0 255] og [-1000 10000] og ændre på parametre i parameter.txt file men det er lort uanset hvad."
$run2 = " Tror jeg prøver med et større volumen nu.."

$xml5 = "<w:p xmlns:w='$wNs'><w:r><w:t>" + $run1 + "</w:t></w:r><w:r><w:t xml:space='preserve'>" + $run2 + "</w:t></w:r></w:p>"
$p5.Range.InsertXML($xml5)

Write-Output "Inserted paragraphs 4 and 5 (count now $($d.Paragraphs.Count))"
